$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Find the last used row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# Insert 5 new columns right before column BJ (the last column, which holds
# the per-row "count" value). This shifts the existing last column 5 places
# to the right (BJ -> BO) and creates 5 blank columns in its place.
$ws.Columns("BJ:BN").Insert()

# The 5 newly inserted columns (BJ:BN) should be filled with the same
# "group" value already present on the row (mirroring columns C:BI), for
# every data row below the header.
for ($r = 2; $r -le $lastRow; $r++) {
    $groupValue = $ws.Range("BI" + $r).Value2
    $ws.Range("BJ" + $r + ":BN" + $r).Value2 = $groupValue
}
